$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the current row 5 ("Page d'accueil" / Home page),
# pushing the Sprint 2 / Sprint 3 rows down by two.
$ws.Rows("6:7").Insert()

# Copy the formatting from an existing body row pair so the new rows pick up
# the same table-style borders/fills as their neighbours (plain Insert leaves
# them unbordered).
$ws.Range("A2:E3").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new rows with the "propriétaire" (owner) account stories.
$ws.Range("A6").Value = "Sprint 1"
$ws.Range("B6").Value = "Créer un compte propriétaire"
$ws.Range("C6").Value = "Haute"
$ws.Range("D6").Value = "À faire"
$ws.Range("E6").Value = "Formulaire inscription"

$ws.Range("A7").Value = "Sprint 1"
$ws.Range("B7").Value = "Se connecter Propriétaire"
$ws.Range("C7").Value = "Haute"
$ws.Range("D7").Value = "À faire"
$ws.Range("E7").Value = "Page login"

$ws.Range("E21").Select() | Out-Null
